$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: all columns A-R get 99.380798999990645, except F2 = 546.59439449994863
$row2Value = 99.380798999990645
$row2F = 546.59439449994863

# Row 3: all columns A-R get 88.888888888888886, except F3 = 488.88888888888891
$row3Value = 88.888888888888886
$row3F = 488.88888888888891

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

foreach ($col in $cols) {
    if ($col -eq "F") {
        $ws.Range("F2").Value = $row2F
        $ws.Range("F3").Value = $row3F
    } else {
        $ws.Range("$col" + "2").Value = $row2Value
        $ws.Range("$col" + "3").Value = $row3Value
    }
}
